# Weekly update: add a new week of "Acelga" (Vega Central Mapocho de Santiago) data.
# This inserts two new data rows (440 and 441) above the existing row 440, which
# pushes all subsequent rows (old 440..520) down by two positions (new 442..522).
# The two new rows carry the new week's "Primera" and "Segunda" quality records.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new blank rows at position 440; existing rows 440+ shift down to 442+.
$ws.Rows.Item(440).Resize(2).Insert()

# ---- Row 440: new "Primera" record for 2022-03-17 (serial 44637) ----
$ws.Range("A440").Value = 9
$ws.Range("B440").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C440").Value = 'Metropolitana'
$ws.Range("D440").Value = 44637
$ws.Range("E440").Value = 13
$ws.Range("F440").Value = 100112009
$ws.Range("G440").Value = 'Acelga'
$ws.Range("H440").Value = 'Sin especificar'
$ws.Range("I440").Value = 'Primera'
$ws.Range("J440").Value = 61
$ws.Range("K440").Value = 18000
$ws.Range("L440").Value = 18000
$ws.Range("M440").Value = 18000
$ws.Range("N440").Value = '$/docena de atados'
$ws.Range("O440").Value = 'Región Metropolitana'
$ws.Range("P440").Value = 6000
$ws.Range("Q440").Value = 3
$ws.Range("R440").Value = 'Hortaliza'

# ---- Row 441: new "Segunda" record for 2022-03-17 (serial 44637) ----
$ws.Range("A441").Value = 9
$ws.Range("B441").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C441").Value = 'Metropolitana'
$ws.Range("D441").Value = 44637
$ws.Range("E441").Value = 13
$ws.Range("F441").Value = 100112009
$ws.Range("G441").Value = 'Acelga'
$ws.Range("H441").Value = 'Sin especificar'
$ws.Range("I441").Value = 'Segunda'
$ws.Range("J441").Value = 43
$ws.Range("K441").Value = 16000
$ws.Range("L441").Value = 16000
$ws.Range("M441").Value = 16000
$ws.Range("N441").Value = '$/docena de atados'
$ws.Range("O441").Value = 'Región Metropolitana'
$ws.Range("P441").Value = 5333
$ws.Range("Q441").Value = 3
$ws.Range("R441").Value = 'Hortaliza'
